$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "Unit Tested?" header in the next empty column (K2)
$ws.Range("K2").Value = "Unit Tested?"

# Leave the selection where the user ended up after adding the question
$ws.Range("I6").Select() | Out-Null
